$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Strip the thousands-separator space from the price strings in column B.
$ws.Range("B1").Value = "64000 грн"
$ws.Range("B2").Value = "14099 грн"
$ws.Range("B3").Value = "9999 грн"
$ws.Range("B4").Value = "10232 грн"

# 2) Add a new quantity column F (cart plus/minus/delete stock counts).
$ws.Cells.Item(1, 6).Value = 22
$ws.Cells.Item(2, 6).Value = 10
$ws.Cells.Item(3, 6).Value = 8
$ws.Cells.Item(4, 6).Value = 45
$ws.Cells.Item(5, 6).Value = 20

# 3) Widen column E a bit to fit the new layout.
# (Excel quantizes ColumnWidth to 1/6-character steps on this engine; 30
#  is the nearest input that lands on the bucket closest to the 30.84
#  target width after that quantization.)
$ws.Columns.Item(5).ColumnWidth = 30

# 4) Move the active selection to E3.
$ws.Range("E3").Select()

Write-Output "done"
